$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = "[Aderci-Tornearia-2A, Andre B.-Elet. Dig. Bas.-2A, Aline S. M.-Metalografia-2A, Anderson-Ajustagem-2A]"

$ws.Range("B11").Value = "-"
$ws.Range("E11").Value = "[Claudinei-Des. Maq. Cad_T2-2A, Claudinei-Des. Maq. Cad_T2-2A]"
$ws.Range("F11").Value = "Ludoff-Máquinas Térmicas e de Fluxo"

$ws.Range("B12").Value = "-"
$ws.Range("E12").Value = "[Claudinei-Des. Maq. Cad_T2-2A, Suzanny-Des. Maq. Cad_T1-2A]"
$ws.Range("F12").Value = "Ludoff-Máquinas Térmicas e de Fluxo"

$ws.Range("B14").Value = "-"
$ws.Range("E14").Value = "[Suzanny-Des. Maq. Cad_T1-2A, Suzanny-Des. Maq. Cad_T1-2A]"
$ws.Range("F14").Value = "[Anderson-Ajustagem-2A, Aline S. M.-Metalografia-2A, Aderci-Tornearia-2A, Andre B.-Elet. Dig. Bas.-2A]"

$ws.Range("B15").Value = "-"
$ws.Range("E15").Value = "Gilberto-Mec. Tec. Res. Mat."
$ws.Range("F15").Value = "[Anderson-Ajustagem-2A, Aline S. M.-Metalografia-2A, Aderci-Tornearia-2A, Andre B.-Elet. Dig. Bas.-2A]"

$ws.Range("B16").Value = "-"
$ws.Range("E16").Value = "Gilberto-Mec. Tec. Res. Mat."
$ws.Range("F16").Value = "[Anderson-Ajustagem-2A, Aline S. M.-Metalografia-2A, Aderci-Tornearia-2A, Andre B.-Elet. Dig. Bas.-2A]"
